$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every touched D/E (and occasionally B/C) cell holds plain text in the
# source workbook (inline strings, t="inlineStr") -- including values
# that look numeric, like "301.95" or the multi-dot "thousands" prices
# such as "42.990.99", plus percentage deltas with leading/trailing
# padding spaces such as "  -0.38%  ". Assigning straight to .Value lets
# Excel auto-convert anything that parses as a number/date, which is not
# what the source file stores. So for each edited row we flip the
# contiguous run of touched columns to the "@" (Text) number format
# first, assign the literal strings, then flip the same range's .Style
# back to "Normal" so the saved cells end up with no explicit style
# index again (matching the original, unstyled D/E cells). Doing this
# per contiguous row-range (instead of per cell) means the workbook ends
# up with a single reused throwaway "Text" style-table entry rather than
# one per cell.

$rowRange = $ws.Range('D2:E2')
$rowRange.NumberFormat = "@"
$ws.Range('D2').Value = '42.990.99'
$ws.Range('E2').Value = '  -0.38%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D3:E3')
$rowRange.NumberFormat = "@"
$ws.Range('D3').Value = '2.302.14'
$ws.Range('E3').Value = '  -0.61%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E4')
$rowRange.NumberFormat = "@"
$ws.Range('E4').Value = '  +0.03%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D5:E5')
$rowRange.NumberFormat = "@"
$ws.Range('D5').Value = '301.95'
$ws.Range('E5').Value = '  -0.46%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D6:E6')
$rowRange.NumberFormat = "@"
$ws.Range('D6').Value = '98.47'
$ws.Range('E6').Value = '  -3.56%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D7:E7')
$rowRange.NumberFormat = "@"
$ws.Range('D7').Value = '0.524'
$ws.Range('E7').Value = '  +3.59%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E8')
$rowRange.NumberFormat = "@"
$ws.Range('E8').Value = '  +0.01%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E9')
$rowRange.NumberFormat = "@"
$ws.Range('E9').Value = '  +0.44%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D10:E10')
$rowRange.NumberFormat = "@"
$ws.Range('D10').Value = '35.62'
$ws.Range('E10').Value = '  -0.95%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E11')
$rowRange.NumberFormat = "@"
$ws.Range('E11').Value = '  -0.86%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E12')
$rowRange.NumberFormat = "@"
$ws.Range('E12').Value = '  -1.11%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D13:E13')
$rowRange.NumberFormat = "@"
$ws.Range('D13').Value = '17.94'
$ws.Range('E13').Value = '  -0.09%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D15:E15')
$rowRange.NumberFormat = "@"
$ws.Range('D15').Value = '2.663.28'
$ws.Range('E15').Value = '  -1.08%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D16:E16')
$rowRange.NumberFormat = "@"
$ws.Range('D16').Value = '2.313.65'
$ws.Range('E16').Value = '  -0.12%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E17')
$rowRange.NumberFormat = "@"
$ws.Range('E17').Value = '  -3.21%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D18:E18')
$rowRange.NumberFormat = "@"
$ws.Range('D18').Value = '42.901.26'
$ws.Range('E18').Value = '  -0.42%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D19:E19')
$rowRange.NumberFormat = "@"
$ws.Range('D19').Value = '13.40'
$ws.Range('E19').Value = '  +5.79%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('B20:E20')
$rowRange.NumberFormat = "@"
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0908'
$ws.Range('E20').Value = '  +0.33%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('B21:E21')
$rowRange.NumberFormat = "@"
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '6.17'
$ws.Range('E21').Value = '  -0.14%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E22')
$rowRange.NumberFormat = "@"
$ws.Range('E22').Value = '  +0.48%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D23:E23')
$rowRange.NumberFormat = "@"
$ws.Range('D23').Value = '239.54'
$ws.Range('E23').Value = '  +0.82%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E24')
$rowRange.NumberFormat = "@"
$ws.Range('E24').Value = '  -2.98%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('B25:E25')
$rowRange.NumberFormat = "@"
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '2.45'
$ws.Range('E25').Value = '  -1.15%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('B26:E26')
$rowRange.NumberFormat = "@"
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.06%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D27:E27')
$rowRange.NumberFormat = "@"
$ws.Range('D27').Value = '24.75'
$ws.Range('E27').Value = '  -0.40%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D28:E28')
$rowRange.NumberFormat = "@"
$ws.Range('D28').Value = '168.04'
$ws.Range('E28').Value = '  -0.10%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D29:E29')
$rowRange.NumberFormat = "@"
$ws.Range('D29').Value = '9.13'
$ws.Range('E29').Value = '  -1.27%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E30')
$rowRange.NumberFormat = "@"
$ws.Range('E30').Value = '  -12.63%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D31:E31')
$rowRange.NumberFormat = "@"
$ws.Range('D31').Value = '33.30'
$ws.Range('E31').Value = '  -3.53%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D32:E32')
$rowRange.NumberFormat = "@"
$ws.Range('D32').Value = '5.19'
$ws.Range('E32').Value = '  +2.91%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D33:E33')
$rowRange.NumberFormat = "@"
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.01%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D34:E34')
$rowRange.NumberFormat = "@"
$ws.Range('D34').Value = '4.83'
$ws.Range('E34').Value = '  +1.75%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D35:E35')
$rowRange.NumberFormat = "@"
$ws.Range('D35').Value = '18.20'
$ws.Range('E35').Value = '  +5.48%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D36:E36')
$rowRange.NumberFormat = "@"
$ws.Range('D36').Value = '2.40'
$ws.Range('E36').Value = '  -0.51%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D37:E37')
$rowRange.NumberFormat = "@"
$ws.Range('D37').Value = '0.0691'
$ws.Range('E37').Value = '  -0.34%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E38')
$rowRange.NumberFormat = "@"
$ws.Range('E38').Value = '  -1.45%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E39')
$rowRange.NumberFormat = "@"
$ws.Range('E39').Value = '  -0.17%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E40')
$rowRange.NumberFormat = "@"
$ws.Range('E40').Value = '  +1.30%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E41')
$rowRange.NumberFormat = "@"
$ws.Range('E41').Value = '  -3.00%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D42:E42')
$rowRange.NumberFormat = "@"
$ws.Range('D42').Value = '1.997.25'
$ws.Range('E42').Value = '  +0.19%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E43')
$rowRange.NumberFormat = "@"
$ws.Range('E43').Value = '  -0.51%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('B44:E44')
$rowRange.NumberFormat = "@"
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').Value = '2.13'
$ws.Range('E44').Value = '  -7.35%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('B45:E45')
$rowRange.NumberFormat = "@"
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '10.07'
$ws.Range('E45').Value = '  -1.85%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D46:E46')
$rowRange.NumberFormat = "@"
$ws.Range('D46').Value = '17.44'
$ws.Range('E46').Value = '  -1.26%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('E47')
$rowRange.NumberFormat = "@"
$ws.Range('E47').Value = '  -2.49%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D48:E48')
$rowRange.NumberFormat = "@"
$ws.Range('D48').Value = '54.70'
$ws.Range('E48').Value = '  -2.44%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('D49:E49')
$rowRange.NumberFormat = "@"
$ws.Range('D49').Value = '2.533.34'
$ws.Range('E49').Value = '  +0.43%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('B50:E50')
$rowRange.NumberFormat = "@"
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '1.53'
$ws.Range('E50').Value = '  -0.22%  '
$rowRange.Style = "Normal"

$rowRange = $ws.Range('B51:E51')
$rowRange.NumberFormat = "@"
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '73.33'
$ws.Range('E51').Value = '  +4.22%  '
$rowRange.Style = "Normal"
